$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns for team record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in team record values for every data row (rows 2-45)
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 91   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 71   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
